$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Old rows 10 and 11 were both "Contact" / "No display for ContactDetail"
# Row 10 becomes Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 becomes Description / Extended class codes for encounters
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Extended class codes for encounters"

# Old row 12 (Description / Extended class codes for encounters) is removed entirely,
# shifting old rows 13,14,15 up to 12,13,14
$ws.Rows("12").Delete()
